$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("N2").Value = 13
$ws.Range("O2").Value = 1.22
$ws.Range("P2").Value = 4
$ws.Range("Q2").Value = 1.75
$ws.Range("R2").Value = 2.05
$ws.Range("AD2").Value = 8.5
$ws.Range("AO2").Value = 7.5

# Row 4 updates
$ws.Range("G4").Value = 1.85
$ws.Range("I4").Value = 4
$ws.Range("Q4").Value = 2.15
$ws.Range("R4").Value = 1.67
$ws.Range("U4").Value = 1.91
$ws.Range("V4").Value = 1.8
$ws.Range("W4").Value = 6.5
$ws.Range("X4").Value = 8.5
$ws.Range("AB4").Value = 29
$ws.Range("AF4").Value = 51
$ws.Range("AL4").Value = 34
$ws.Range("AO4").Value = 11
$ws.Range("AQ4").Value = 41
$ws.Range("AU4").Value = 8.5
$ws.Range("BA4").Value = 101
